$p = $ppt.ActivePresentation
$newDate = "6/16/2023"
$ppPlaceholderDate = 16

# NB: Slides.Item(n).Master.CustomLayouts.Item(k) mis-resolves to the
# *first* layout for every k in this host (reads and writes alike), so
# custom layouts are reached via Presentation.Designs.Item(1).SlideMaster
# instead, which resolves each layout correctly.
$design = $p.Designs.Item(1)
$master = $design.SlideMaster

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $isDatePlaceholder = $false
            try {
                if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                    $isDatePlaceholder = $true
                }
            } catch {
                $isDatePlaceholder = $false
            }
            if ((-not $isDatePlaceholder) -and ($shp.Name -like "Date Placeholder*")) {
                $isDatePlaceholder = $true
            }
            if ($isDatePlaceholder) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# The slide master's own "Date Placeholder" (datetimeFigureOut field).
Update-DatePlaceholder $master.Shapes

# Every slide layout's "Date Placeholder".
for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    Update-DatePlaceholder $layout.Shapes
}
